$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("I1:I233")
$range.Replace("http://pixmosaic.ru", "https://pixmosaic.ru")
